$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price, Volume(1h)) updates, per commit "Updated cryptos list"
$updates = @(
    @{ Row = 2; Price = '63.702.73'; Volume = '  -0.93%  ' }
    @{ Row = 3; Price = '3.421.98'; Volume = '  -2.36%  ' }
    @{ Row = 4; Price = '0.999'; Volume = '  -0.04%  ' }
    @{ Row = 5; Price = '579.27'; Volume = '  -1.87%  ' }
    @{ Row = 6; Price = '128.94'; Volume = '  -4.07%  ' }
    @{ Row = 7; Price = '1.00'; Volume = '  +0.01%  ' }
    @{ Row = 8; Price = '0.480'; Volume = '  -1.58%  ' }
    @{ Row = 9; Price = '7.56'; Volume = '  +3.36%  ' }
    @{ Row = 10; Price = '0.124'; Volume = '  -0.17%  ' }
    @{ Row = 11; Price = '0.381'; Volume = '  -1.24%  ' }
    @{ Row = 12; Price = '4.000.98'; Volume = '  -2.41%  ' }
    @{ Row = 13; Price = '0.119'; Volume = '  -0.43%  ' }
    @{ Row = 14; Price = '0.0000176'; Volume = '  -2.73%  ' }
    @{ Row = 15; Price = '3.424.61'; Volume = '  -2.31%  ' }
    @{ Row = 16; Price = '63.676.21'; Volume = '  -0.99%  ' }
    @{ Row = 17; Price = '25.28'; Volume = '  -1.49%  ' }
    @{ Row = 18; Price = '9.81'; Volume = '  -0.57%  ' }
    @{ Row = 19; Price = '5.63'; Volume = '  -2.02%  ' }
    @{ Row = 20; Price = '13.31'; Volume = '  -1.52%  ' }
    @{ Row = 21; Price = '382.69'; Volume = '  -2.72%  ' }
    @{ Row = 22; Price = '0.562'; Volume = '  -1.65%  ' }
    @{ Row = 23; Price = '3.556.35'; Volume = '  -2.38%  ' }
    @{ Row = 24; Price = '73.99'; Volume = '  -0.88%  ' }
    @{ Row = 25; Price = '1.00'; Volume = '  -0.11%  ' }
    @{ Row = 26; Price = '0.0000109'; Volume = '  -5.08%  ' }
    @{ Row = 27; Price = '1.00'; Volume = '  +0.10%  ' }
    @{ Row = 28; Price = '2.19'; Volume = '  -2.94%  ' }
    @{ Row = 29; Price = '7.01'; Volume = '  -4.98%  ' }
    @{ Row = 30; Price = '7.87'; Volume = '  -4.39%  ' }
    @{ Row = 31; Price = '0.153'; Volume = '  -0.76%  ' }
    @{ Row = 32; Price = '1.41'; Volume = '  -4.74%  ' }
    @{ Row = 33; Price = '3.449.85'; Volume = '  -2.15%  ' }
    @{ Row = 34; Price = '1.00'; Volume = '  -0.07%  ' }
    @{ Row = 35; Price = '22.73'; Volume = '  -3.21%  ' }
    @{ Row = 36; Price = '5.13'; Volume = '  -0.14%  ' }
    @{ Row = 37; Price = '6.72'; Volume = '  -2.45%  ' }
    @{ Row = 38; Price = '164.06'; Volume = '  -2.14%  ' }
    @{ Row = 39; Price = '1.51'; Volume = '  -2.61%  ' }
    @{ Row = 40; Price = '0.0768'; Volume = '  -1.56%  ' }
    @{ Row = 41; Price = '0.784'; Volume = '  -3.41%  ' }
    @{ Row = 42; Price = '1.00'; Volume = '  -0.04%  ' }
    @{ Row = 43; Price = '41.42'; Volume = '  -0.90%  ' }
    @{ Row = 44; Price = '4.29'; Volume = '  -2.41%  ' }
    @{ Row = 45; Price = '1.60'; Volume = '  -3.68%  ' }
    @{ Row = 46; Price = '23.15'; Volume = '  -7.50%  ' }
    @{ Row = 47; Price = '1.10'; Volume = '  -6.02%  ' }
    @{ Row = 48; Price = '6.69'; Volume = '  -0.92%  ' }
    @{ Row = 49; Price = '0.886'; Volume = '  -0.83%  ' }
    @{ Row = 50; Price = '2.278.83'; Volume = '  -3.27%  ' }
    @{ Row = 51; Price = '0.0252'; Volume = '  -2.66%  ' }
)

foreach ($u in $updates) {
    $priceCell = $ws.Range("D" + $u.Row)
    # Many prices look like plain numbers (e.g. "1.00", "0.999"); force them to
    # stay text so trailing zeros / precision survive exactly as scraped, then
    # restore the default (unstyled) cell style so no formatting is introduced.
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u.Price
    $priceCell.Style = "Normal"

    $ws.Range("E" + $u.Row).Value = $u.Volume
}
